$d = $word.ActiveDocument

# --- Change 1: split "进入项目子目录" into "进入项目" + "根" + "目录" ---
$found = $d.Content.Find.Execute("进入项目子目录", $true, $false, $false, $false, $false,
                                  $true, 1, $false, "进入项目根目录", 2)
